$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 2308
$ws.Range("I127").Value = 790.6667
$ws.Range("J127").Value = 3066.6667
$ws.Range("K127").Value = 2372.0001
$ws.Range("L127").Value = 9200.000100000001
$ws.Range("M127").Value = 2587.9999
$ws.Range("N127").Value = -19120.0001

$ws.Range("H132").Value = 3103.4
$ws.Range("I132").Value = 2975
$ws.Range("J132").Value = 3403
$ws.Range("K132").Value = 8925
$ws.Range("L132").Value = 10209
$ws.Range("M132").Value = -6395
$ws.Range("N132").Value = -15269

$ws.Range("H133").Value = 95796.664
$ws.Range("J133").Value = 95796.664
$ws.Range("L133").Value = 95796.664
$ws.Range("N133").Value = -105916.664

$ws.Range("H137").Value = 2077.5588
$ws.Range("I137").Value = 1475.4667
$ws.Range("J137").Value = 2552.8948
$ws.Range("K137").Value = 4426.4001
$ws.Range("L137").Value = 7658.6844
$ws.Range("M137").Value = -1876.4001
$ws.Range("N137").Value = -12758.6844

$ws.Range("H141").Value = 6483.3335
$ws.Range("I141").Value = 1823.6842
$ws.Range("J141").Value = 50750
$ws.Range("K141").Value = 5471.0526
$ws.Range("L141").Value = 152250
$ws.Range("M141").Value = -291.0526
$ws.Range("N141").Value = -162610

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 47500
$ws.Range("I9").Value = 50000
$ws.Range("J9").Value = 45000
$ws.Range("K9").Value = 50000
$ws.Range("L9").Value = 45000
$ws.Range("M9").Value = -49830
$ws.Range("N9").Value = -45340

$ws.Range("H18").Value = 32650
$ws.Range("I18").Value = 25300
$ws.Range("J18").Value = 40000
$ws.Range("K18").Value = 25300
$ws.Range("L18").Value = 40000
$ws.Range("M18").Value = -24978
$ws.Range("N18").Value = -40644

$ws.Range("H20").Value = 47500
$ws.Range("I20").Value = 50000
$ws.Range("J20").Value = 45000
$ws.Range("K20").Value = 50000
$ws.Range("L20").Value = 45000
$ws.Range("M20").Value = -49730
$ws.Range("N20").Value = -45540

$ws.Range("H45").Value = 2181
$ws.Range("I45").Value = 1824.6666
$ws.Range("K45").Value = 1824.6666
$ws.Range("M45").Value = -1447.6666

$ws.Range("H132").Value = 8984.833000000001
$ws.Range("I132").Value = 10824.23
$ws.Range("K132").Value = 32472.69
$ws.Range("M132").Value = -29942.69

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 84843.664
$ws.Range("I20").Value = 101232.4
$ws.Range("J20").Value = 2900
$ws.Range("K20").Value = 101232.4
$ws.Range("L20").Value = 2900
$ws.Range("M20").Value = -100985.4
$ws.Range("N20").Value = -3394

$ws.Range("H134").Value = 2566.375
$ws.Range("I134").Value = 2214.2942
$ws.Range("K134").Value = 6642.882599999999
$ws.Range("M134").Value = -4107.882599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2345.3
$ws.Range("J31").Value = 2793.5715
$ws.Range("L31").Value = 2793.5715
$ws.Range("N31").Value = -3383.5715

$ws.Range("H34").Value = 2345.3
$ws.Range("J34").Value = 2793.5715
$ws.Range("L34").Value = 2793.5715
$ws.Range("N34").Value = -3197.5715

$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372

$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 167582.53
$ws.Range("I68").Value = 333953.56
$ws.Range("J68").Value = 1211.5
$ws.Range("K68").Value = 1001860.68
$ws.Range("L68").Value = 3634.5
$ws.Range("M68").Value = -1001049.68
$ws.Range("N68").Value = -5256.5

$ws.Range("H71").Value = 167582.53
$ws.Range("I71").Value = 333953.56
$ws.Range("J71").Value = 1211.5
$ws.Range("K71").Value = 3005582.04
$ws.Range("L71").Value = 10903.5
$ws.Range("M71").Value = -3001526.04
$ws.Range("N71").Value = -19015.5

$ws.Range("H107").Value = 1063.5902
$ws.Range("I107").Value = 715.89795
$ws.Range("J107").Value = 2483.3333
$ws.Range("K107").Value = 2147.69385
$ws.Range("L107").Value = 7449.999899999999
$ws.Range("M107").Value = -227.6938500000001
$ws.Range("N107").Value = -11289.9999

$ws.Range("H113").Value = 303702.16
$ws.Range("I113").Value = 667219.5600000001
$ws.Range("J113").Value = 771
$ws.Range("K113").Value = 2001658.68
$ws.Range("L113").Value = 2313
$ws.Range("M113").Value = -1999488.68
$ws.Range("N113").Value = -6653

$ws.Range("H137").Value = 3242.3
$ws.Range("I137").Value = 2052.8572
$ws.Range("J137").Value = 6017.6665
$ws.Range("K137").Value = 6158.571599999999
$ws.Range("L137").Value = 18052.9995
$ws.Range("M137").Value = -1058.571599999999
$ws.Range("N137").Value = -28252.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5333.3335
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 5333.3335
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 5333.3335
$ws.Range("N43").Value = -5635.3335
$ws.Range("M43").ClearContents()

$ws.Range("H93").Value = 31642.857
$ws.Range("J93").Value = 31642.857
$ws.Range("L93").Value = 31642.857
$ws.Range("N93").Value = -35386.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1830.7727
$ws.Range("I136").Value = 1328.1177
$ws.Range("J136").Value = 3539.8
$ws.Range("K136").Value = 3984.3531
$ws.Range("L136").Value = 10619.4
$ws.Range("M136").Value = -1434.3531
$ws.Range("N136").Value = -15719.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 42400
$ws.Range("J8").Value = 42400
$ws.Range("L8").Value = 42400
$ws.Range("N8").Value = -42680

$ws.Range("H38").Value = 7030.6665
$ws.Range("I38").Value = 546
$ws.Range("J38").Value = 20000
$ws.Range("K38").Value = 546
$ws.Range("L38").Value = 20000
$ws.Range("M38").Value = -73
$ws.Range("N38").Value = -20946

$ws.Range("H39").Value = 6585523.5
$ws.Range("I39").Value = 13131380
$ws.Range("J39").Value = 39666.668
$ws.Range("K39").Value = 13131380
$ws.Range("L39").Value = 39666.668
$ws.Range("M39").Value = -13130967
$ws.Range("N39").Value = -40492.668

$ws.Range("H49").Value = 50000
$ws.Range("J49").Value = 50000
$ws.Range("L49").Value = 50000
$ws.Range("N49").Value = -50460

$ws.Range("H136").Value = 3321
$ws.Range("I136").Value = 3924.25
$ws.Range("J136").Value = 2882.2727
$ws.Range("K136").Value = 11772.75
$ws.Range("L136").Value = 8646.8181
$ws.Range("M136").Value = -9222.75
$ws.Range("N136").Value = -13746.8181
Write-Host "Edit complete"
